# Update scraped_at timestamps in the "snapshot" sheet (K2:K52).
# This mirrors a re-run of the KHL injuries scraper ~4 hours later,
# refreshing the per-row scrape timestamp for each injured player entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-10T07:01:26.910693+00:00"
$ws.Range("K3").Value = "2025-11-10T07:01:26.910729+00:00"
$ws.Range("K4").Value = "2025-11-10T07:01:26.910750+00:00"
$ws.Range("K5").Value = "2025-11-10T07:01:29.329086+00:00"
$ws.Range("K6").Value = "2025-11-10T07:01:29.329118+00:00"
$ws.Range("K7").Value = "2025-11-10T07:01:29.329139+00:00"
$ws.Range("K8").Value = "2025-11-10T07:01:31.645250+00:00"
$ws.Range("K9").Value = "2025-11-10T07:01:34.014152+00:00"
$ws.Range("K10").Value = "2025-11-10T07:01:34.014182+00:00"
$ws.Range("K11").Value = "2025-11-10T07:01:34.014201+00:00"
$ws.Range("K12").Value = "2025-11-10T07:01:36.803766+00:00"
$ws.Range("K13").Value = "2025-11-10T07:01:36.803796+00:00"
$ws.Range("K14").Value = "2025-11-10T07:01:36.803815+00:00"
$ws.Range("K15").Value = "2025-11-10T07:01:36.803831+00:00"
$ws.Range("K16").Value = "2025-11-10T07:01:42.335403+00:00"
$ws.Range("K17").Value = "2025-11-10T07:01:44.793298+00:00"
$ws.Range("K18").Value = "2025-11-10T07:01:47.105772+00:00"
$ws.Range("K19").Value = "2025-11-10T07:01:47.105801+00:00"
$ws.Range("K20").Value = "2025-11-10T07:01:47.105819+00:00"
$ws.Range("K21").Value = "2025-11-10T07:01:49.451508+00:00"
$ws.Range("K22").Value = "2025-11-10T07:01:52.210226+00:00"
$ws.Range("K23").Value = "2025-11-10T07:01:52.210259+00:00"
$ws.Range("K24").Value = "2025-11-10T07:01:54.980169+00:00"
$ws.Range("K25").Value = "2025-11-10T07:01:54.980205+00:00"
$ws.Range("K26").Value = "2025-11-10T07:01:54.980230+00:00"
$ws.Range("K27").Value = "2025-11-10T07:01:57.292229+00:00"
$ws.Range("K28").Value = "2025-11-10T07:01:57.292258+00:00"
$ws.Range("K29").Value = "2025-11-10T07:01:57.292277+00:00"
$ws.Range("K30").Value = "2025-11-10T07:01:57.292296+00:00"
$ws.Range("K31").Value = "2025-11-10T07:01:57.292311+00:00"
$ws.Range("K32").Value = "2025-11-10T07:01:59.997505+00:00"
$ws.Range("K33").Value = "2025-11-10T07:01:59.997538+00:00"
$ws.Range("K34").Value = "2025-11-10T07:02:02.276786+00:00"
$ws.Range("K35").Value = "2025-11-10T07:02:02.276816+00:00"
$ws.Range("K36").Value = "2025-11-10T07:02:02.276841+00:00"
$ws.Range("K37").Value = "2025-11-10T07:02:04.640557+00:00"
$ws.Range("K38").Value = "2025-11-10T07:02:04.640586+00:00"
$ws.Range("K39").Value = "2025-11-10T07:02:04.640605+00:00"
$ws.Range("K40").Value = "2025-11-10T07:02:07.071729+00:00"
$ws.Range("K41").Value = "2025-11-10T07:02:07.071766+00:00"
$ws.Range("K42").Value = "2025-11-10T07:02:07.071786+00:00"
$ws.Range("K43").Value = "2025-11-10T07:02:07.071805+00:00"
$ws.Range("K44").Value = "2025-11-10T07:02:07.071823+00:00"
$ws.Range("K45").Value = "2025-11-10T07:02:07.071839+00:00"
$ws.Range("K46").Value = "2025-11-10T07:02:09.304379+00:00"
$ws.Range("K47").Value = "2025-11-10T07:02:09.304409+00:00"
$ws.Range("K48").Value = "2025-11-10T07:02:14.138463+00:00"
$ws.Range("K49").Value = "2025-11-10T07:02:14.138505+00:00"
$ws.Range("K50").Value = "2025-11-10T07:02:14.138528+00:00"
$ws.Range("K51").Value = "2025-11-10T07:02:16.453671+00:00"
$ws.Range("K52").Value = "2025-11-10T07:02:16.453702+00:00"
